# Auto-generated script to update TPM-derived NATMI metrics in App-Cd74 sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 89.38217433333334
$ws.Range("H2").Value = 268.146523
$ws.Range("I2").Value = 0.2143552015363441
$ws.Range("J2").Value = 0.2175965347165783
$ws.Range("M2").Value = 12.492041
$ws.Range("N2").Value = 37.476123
$ws.Range("O2").Value = 0.001845183434243001
$ws.Range("P2").Value = 0.001846965528678714
$ws.Range("Q2").Value = 1116.565786441148
$ws.Range("R2").Value = 10049.09207797033
$ws.Range("S2").Value = 0.0003955246669186819
$ws.Range("T2").Value = 0.0004018932987814612
$ws.Range("G3").Value = 89.38217433333334
$ws.Range("H3").Value = 268.146523
$ws.Range("I3").Value = 0.2143552015363441
$ws.Range("J3").Value = 0.2175965347165783
$ws.Range("O3").Value = 0.0009361563262120847
$ws.Range("P3").Value = 0.0009370604742490439
$ws.Range("Q3").Value = 566.491171127214
$ws.Range("R3").Value = 5098.420540144926
$ws.Range("S3").Value = 0.0002006699779747149
$ws.Range("T3").Value = 0.0002039011120164654
$ws.Range("G4").Value = 89.38217433333334
$ws.Range("H4").Value = 268.146523
$ws.Range("I4").Value = 0.2143552015363441
$ws.Range("J4").Value = 0.2175965347165783
$ws.Range("M4").Value = 4688.500488333333
$ws.Range("N4").Value = 14065.501465
$ws.Range("O4").Value = 0.6925324238459419
$ws.Range("P4").Value = 0.6932012777691797
$ws.Range("Q4").Value = 419068.3680101285
$ws.Range("R4").Value = 3771615.312091156
$ws.Range("S4").Value = 0.1484479272839498
$ws.Range("T4").Value = 0.1508381959036778
$ws.Range("G5").Value = 89.38217433333334
$ws.Range("H5").Value = 268.146523
$ws.Range("I5").Value = 0.2143552015363441
$ws.Range("J5").Value = 0.2175965347165783
$ws.Range("M5").Value = 19.5968845
$ws.Range("N5").Value = 39.193769
$ws.Range("O5").Value = 0.002894630800697294
$ws.Range("P5").Value = 0.001931617640437256
$ws.Range("Q5").Value = 1751.612146769198
$ws.Range("R5").Value = 10509.67288061519
$ws.Range("S5").Value = 0.0006204791686567775
$ws.Range("T5").Value = 0.0004203133049565604
$ws.Range("G6").Value = 89.38217433333334
$ws.Range("H6").Value = 268.146523
$ws.Range("I6").Value = 0.2143552015363441
$ws.Range("J6").Value = 0.2175965347165783
$ws.Range("M6").Value = 2043.153564333333
$ws.Range("N6").Value = 6129.460693
$ws.Range("O6").Value = 0.3017916055929057
$ws.Range("P6").Value = 0.3020830785874553
$ws.Range("Q6").Value = 182621.5080770134
$ws.Range("R6").Value = 1643593.57269312
$ws.Range("S6").Value = 0.06469060043884417
$ws.Range("T6").Value = 0.06573223109714607
$ws.Range("I7").Value = 0.2934277926151677
$ws.Range("J7").Value = 0.2978648075949286
$ws.Range("M7").Value = 12.492041
$ws.Range("N7").Value = 37.476123
$ws.Range("O7").Value = 0.001845183434243001
$ws.Range("P7").Value = 0.001846965528678714
$ws.Range("Q7").Value = 1528.451055429576
$ws.Range("R7").Value = 13756.05949886619
$ws.Range("S7").Value = 0.000541428102079998
$ws.Range("T7").Value = 0.0005501460318343505
$ws.Range("I8").Value = 0.2934277926151677
$ws.Range("J8").Value = 0.2978648075949286
$ws.Range("O8").Value = 0.0009361563262120847
$ws.Range("P8").Value = 0.0009370604742490439
$ws.Range("S8").Value = 0.0002746942843431368
$ws.Range("T8").Value = 0.000279117337867004
$ws.Range("I9").Value = 0.2934277926151677
$ws.Range("J9").Value = 0.2978648075949286
$ws.Range("M9").Value = 4688.500488333333
$ws.Range("N9").Value = 14065.501465
$ws.Range("O9").Value = 0.6925324238459419
$ws.Range("P9").Value = 0.6932012777691797
$ws.Range("Q9").Value = 573656.7403016982
$ws.Range("R9").Value = 5162910.662715284
$ws.Range("S9").Value = 0.2032082604435465
$ws.Range("T9").Value = 0.2064802652272753
$ws.Range("I10").Value = 0.2934277926151677
$ws.Range("J10").Value = 0.2978648075949286
$ws.Range("M10").Value = 19.5968845
$ws.Range("N10").Value = 39.193769
$ws.Range("O10").Value = 0.002894630800697294
$ws.Range("P10").Value = 0.001931617640437256
$ws.Range("Q10").Value = 2397.75700361186
$ws.Range("R10").Value = 14386.54202167116
$ws.Range("S10").Value = 0.0008493651262844821
$ws.Range("T10").Value = 0.0005753609168158131
$ws.Range("I11").Value = 0.2934277926151677
$ws.Range("J11").Value = 0.2978648075949286
$ws.Range("M11").Value = 2043.153564333333
$ws.Range("N11").Value = 6129.460693
$ws.Range("O11").Value = 0.3017916055929057
$ws.Range("P11").Value = 0.3020830785874553
$ws.Range("Q11").Value = 249987.9900978538
$ws.Range("R11").Value = 2249891.910880684
$ws.Range("S11").Value = 0.0885540446589136
$ws.Range("T11").Value = 0.08997991808113605
$ws.Range("G12").Value = 90.33462533333334
$ws.Range("H12").Value = 271.003876
$ws.Range("I12").Value = 0.2166393574945233
$ws.Range("J12").Value = 0.2199152301234996
$ws.Range("M12").Value = 12.492041
$ws.Range("N12").Value = 37.476123
$ws.Range("O12").Value = 0.001845183434243001
$ws.Range("P12").Value = 0.001846965528678714
$ws.Range("Q12").Value = 1128.463843383639
$ws.Range("R12").Value = 10156.17459045275
$ws.Range("S12").Value = 0.0003997393536539416
$ws.Range("T12").Value = 0.0004061758492695504
$ws.Range("G13").Value = 90.33462533333334
$ws.Range("H13").Value = 271.003876
$ws.Range("I13").Value = 0.2166393574945233
$ws.Range("J13").Value = 0.2199152301234996
$ws.Range("O13").Value = 0.0009361563262120847
$ws.Range("P13").Value = 0.0009370604742490439
$ws.Range("Q13").Value = 572.5276665073681
$ws.Range("R13").Value = 5152.748998566311
$ws.Range("S13").Value = 0.0002028083050250194
$ws.Range("T13").Value = 0.0002060738698341142
$ws.Range("G14").Value = 90.33462533333334
$ws.Range("H14").Value = 271.003876
$ws.Range("I14").Value = 0.2166393574945233
$ws.Range("J14").Value = 0.2199152301234996
$ws.Range("M14").Value = 4688.500488333333
$ws.Range("N14").Value = 14065.501465
$ws.Range("O14").Value = 0.6925324238459419
$ws.Range("P14").Value = 0.6932012777691797
$ws.Range("Q14").Value = 423533.9349887421
$ws.Range("R14").Value = 3811805.414898678
$ws.Range("S14").Value = 0.1500297793461097
$ws.Range("T14").Value = 0.1524455185225131
$ws.Range("G15").Value = 90.33462533333334
$ws.Range("H15").Value = 271.003876
$ws.Range("I15").Value = 0.2166393574945233
$ws.Range("J15").Value = 0.2199152301234996
$ws.Range("M15").Value = 19.5968845
$ws.Range("N15").Value = 39.193769
$ws.Range("O15").Value = 0.002894630800697294
$ws.Range("P15").Value = 0.001931617640437256
$ws.Range("Q15").Value = 1770.277219008107
$ws.Range("R15").Value = 10621.66331404864
$ws.Range("S15").Value = 0.0006270909568469191
$ws.Range("T15").Value = 0.0004247921379073704
$ws.Range("G16").Value = 90.33462533333334
$ws.Range("H16").Value = 271.003876
$ws.Range("I16").Value = 0.2166393574945233
$ws.Range("J16").Value = 0.2199152301234996
$ws.Range("M16").Value = 2043.153564333333
$ws.Range("N16").Value = 6129.460693
$ws.Range("O16").Value = 0.3017916055929057
$ws.Range("P16").Value = 0.3020830785874553
$ws.Range("Q16").Value = 184567.5117325162
$ws.Range("R16").Value = 1661107.605592646
$ws.Range("S16").Value = 0.06537993953288766
$ws.Range("T16").Value = 0.06643266974397544
$ws.Range("G17").Value = 18.634161
$ws.Range("H17").Value = 37.268322
$ws.Range("I17").Value = 0.0446882095496985
$ws.Range("J17").Value = 0.03024263611988591
$ws.Range("M17").Value = 12.492041
$ws.Range("N17").Value = 37.476123
$ws.Range("O17").Value = 0.001845183434243001
$ws.Range("P17").Value = 0.001846965528678714
$ws.Range("Q17").Value = 232.778703212601
$ws.Range("R17").Value = 1396.672219275606
$ws.Range("S17").Value = 0.00008245794396708354
$ws.Range("T17").Value = 0.00005585710640980304
$ws.Range("G18").Value = 18.634161
$ws.Range("H18").Value = 37.268322
$ws.Range("I18").Value = 0.0446882095496985
$ws.Range("J18").Value = 0.03024263611988591
$ws.Range("O18").Value = 0.0009361563262120847
$ws.Range("P18").Value = 0.0009370604742490439
$ws.Range("Q18").Value = 118.100591830494
$ws.Range("R18").Value = 708.603550982964
$ws.Range("S18").Value = 0.00004183515007704155
$ws.Range("T18").Value = 0.00002833917894504156
$ws.Range("G19").Value = 18.634161
$ws.Range("H19").Value = 37.268322
$ws.Range("I19").Value = 0.0446882095496985
$ws.Range("J19").Value = 0.03024263611988591
$ws.Range("M19").Value = 4688.500488333333
$ws.Range("N19").Value = 14065.501465
$ws.Range("O19").Value = 0.6925324238459419
$ws.Range("P19").Value = 0.6932012777691797
$ws.Range("Q19").Value = 87366.27294818194
$ws.Range("R19").Value = 524197.6376890917
$ws.Range("S19").Value = 0.03094803407678807
$ws.Range("T19").Value = 0.02096423400141326
$ws.Range("G20").Value = 18.634161
$ws.Range("H20").Value = 37.268322
$ws.Range("I20").Value = 0.0446882095496985
$ws.Range("J20").Value = 0.03024263611988591
$ws.Range("M20").Value = 19.5968845
$ws.Range("N20").Value = 39.193769
$ws.Range("O20").Value = 0.002894630800697294
$ws.Range("P20").Value = 0.001931617640437256
$ws.Range("Q20").Value = 365.1715008714044
$ws.Range("R20").Value = 1460.686003485618
$ws.Range("S20").Value = 0.0001293558677905722
$ws.Range("T20").Value = 0.00005841720942249655
$ws.Range("G21").Value = 18.634161
$ws.Range("H21").Value = 37.268322
$ws.Range("I21").Value = 0.0446882095496985
$ws.Range("J21").Value = 0.03024263611988591
$ws.Range("M21").Value = 2043.153564333333
$ws.Range("N21").Value = 6129.460693
$ws.Range("O21").Value = 0.3017916055929057
$ws.Range("P21").Value = 0.3020830785874553
$ws.Range("Q21").Value = 38072.45246551119
$ws.Range("R21").Value = 228434.7147930671
$ws.Range("S21").Value = 0.01348652651107573
$ws.Range("T21").Value = 0.00913578862369531
$ws.Range("G22").Value = 96.27664699999998
$ws.Range("H22").Value = 288.829941
$ws.Range("I22").Value = 0.2308894388042666
$ws.Range("J22").Value = 0.2343807914451077
$ws.Range("M22").Value = 12.492041
$ws.Range("N22").Value = 37.476123
$ws.Range("O22").Value = 0.001845183434243001
$ws.Range("P22").Value = 0.001846965528678714
$ws.Range("Q22").Value = 1202.691821666527
$ws.Range("R22").Value = 10824.22639499874
$ws.Range("S22").Value = 0.0004260333676232957
$ws.Range("T22").Value = 0.0004328932423835486
$ws.Range("G23").Value = 96.27664699999998
$ws.Range("H23").Value = 288.829941
$ws.Range("I23").Value = 0.2308894388042666
$ws.Range("J23").Value = 0.2343807914451077
$ws.Range("O23").Value = 0.0009361563262120847
$ws.Range("P23").Value = 0.0009370604742490439
$ws.Range("Q23").Value = 610.1873322955379
$ws.Range("R23").Value = 5491.685990659841
$ws.Range("S23").Value = 0.0002161486087921722
$ws.Range("T23").Value = 0.0002196289755864188
$ws.Range("G24").Value = 96.27664699999998
$ws.Range("H24").Value = 288.829941
$ws.Range("I24").Value = 0.2308894388042666
$ws.Range("J24").Value = 0.2343807914451077
$ws.Range("M24").Value = 4688.500488333333
$ws.Range("N24").Value = 14065.501465
$ws.Range("O24").Value = 0.6925324238459419
$ws.Range("P24").Value = 0.6932012777691797
$ws.Range("Q24").Value = 451393.1064745958
$ws.Range("R24").Value = 4062537.958271363
$ws.Range("S24").Value = 0.159898422695548
$ws.Range("T24").Value = 0.1624730641143003
$ws.Range("G25").Value = 96.27664699999998
$ws.Range("H25").Value = 288.829941
$ws.Range("I25").Value = 0.2308894388042666
$ws.Range("J25").Value = 0.2343807914451077
$ws.Range("M25").Value = 19.5968845
$ws.Range("N25").Value = 39.193769
$ws.Range("O25").Value = 0.002894630800697294
$ws.Range("P25").Value = 0.001931617640437256
$ws.Range("Q25").Value = 1886.722331306271
$ws.Range("R25").Value = 11320.33398783763
$ws.Range("S25").Value = 0.000668339681118543
$ws.Range("T25").Value = 0.0004527340713350154
$ws.Range("G26").Value = 96.27664699999998
$ws.Range("H26").Value = 288.829941
$ws.Range("I26").Value = 0.2308894388042666
$ws.Range("J26").Value = 0.2343807914451077
$ws.Range("M26").Value = 2043.153564333333
$ws.Range("N26").Value = 6129.460693
$ws.Range("O26").Value = 0.3017916055929057
$ws.Range("P26").Value = 0.3020830785874553
$ws.Range("Q26").Value = 196707.9744801121
$ws.Range("R26").Value = 1770371.770321009
$ws.Range("S26").Value = 0.06968049445118456
$ws.Range("T26").Value = 0.07080247104150242
